# Auto-generated Excel COM-interop script to apply scheduled-runner price/profit updates
# to the Zodiark_Profits leveling sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1039.9
$ws.Range("J17").Value = 1039.9
$ws.Range("L17").Value = 3119.7
$ws.Range("N17").Value = -3455.7
$ws.Range("H18").Value = 1139.1428
$ws.Range("I18").Value = 1139.1428
$ws.Range("K18").Value = 1139.1428
$ws.Range("M18").Value = -855.1428000000001
$ws.Range("H32").Value = 2991.2273
$ws.Range("I32").Value = 2199
$ws.Range("K32").Value = 2199
$ws.Range("M32").Value = -1873
$ws.Range("H33").Value = 101000140
$ws.Range("I33").Value = 1250124.9
$ws.Range("J33").Value = 500000220
$ws.Range("K33").Value = 1250124.9
$ws.Range("L33").Value = 500000220
$ws.Range("M33").Value = -1249895.9
$ws.Range("N33").Value = -500000678
$ws.Range("H92").Value = 4465267
$ws.Range("I92").Value = 841.5454999999999
$ws.Range("J92").Value = 14287003
$ws.Range("K92").Value = 841.5454999999999
$ws.Range("L92").Value = 14287003
$ws.Range("M92").Value = 406.4545000000001
$ws.Range("N92").Value = -14289499
$ws.Range("H111").Value = 2243.7144
$ws.Range("I111").Value = 1904.4
$ws.Range("K111").Value = 5713.200000000001
$ws.Range("M111").Value = -2646.200000000001
$ws.Range("H113").Value = 5392.923
$ws.Range("J113").Value = 4738.4863
$ws.Range("L113").Value = 4738.4863
$ws.Range("N113").Value = -11246.4863
$ws.Range("H127").Value = 513.5
$ws.Range("I127").Value = 513.5
$ws.Range("K127").Value = 1540.5
$ws.Range("M127").Value = 3419.5
$ws.Range("H132").Value = 1464.15
$ws.Range("I132").Value = 1360.1666
$ws.Range("K132").Value = 4080.4998
$ws.Range("M132").Value = -1550.4998
$ws.Range("H137").Value = 2680.3333
$ws.Range("I137").Value = 2972.7693
$ws.Range("K137").Value = 8918.3079
$ws.Range("M137").Value = -6368.3079
$ws.Range("H138").Value = 2812.946
$ws.Range("I138").Value = 2143
$ws.Range("J138").Value = 2851.2285
$ws.Range("K138").Value = 6429
$ws.Range("L138").Value = 8553.6855
$ws.Range("M138").Value = -1289
$ws.Range("N138").Value = -18833.6855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2185.2727
$ws.Range("I45").Value = 2063.2942
$ws.Range("K45").Value = 2063.2942
$ws.Range("M45").Value = -1686.2942
$ws.Range("H61").Value = 1994.7858
$ws.Range("I61").Value = 1603.7778
$ws.Range("J61").Value = 2698.6
$ws.Range("K61").Value = 1603.7778
$ws.Range("L61").Value = 2698.6
$ws.Range("M61").Value = -1391.7778
$ws.Range("N61").Value = -3122.6
$ws.Range("H74").Value = 3159.7273
$ws.Range("I74").Value = 2935.4546
$ws.Range("K74").Value = 2935.4546
$ws.Range("M74").Value = -2061.4546
$ws.Range("H77").Value = 3159.7273
$ws.Range("I77").Value = 2935.4546
$ws.Range("K77").Value = 14677.273
$ws.Range("M77").Value = -10309.273
$ws.Range("H88").Value = 3367.8235
$ws.Range("J88").Value = 3431.6428
$ws.Range("L88").Value = 3431.6428
$ws.Range("N88").Value = -4243.6428
$ws.Range("H91").Value = 3367.8235
$ws.Range("J91").Value = 3431.6428
$ws.Range("L91").Value = 3431.6428
$ws.Range("N91").Value = -6239.6428
$ws.Range("H97").Value = 592.3125
$ws.Range("I97").Value = 379.72726
$ws.Range("J97").Value = 1060
$ws.Range("K97").Value = 379.72726
$ws.Range("L97").Value = 1060
$ws.Range("M97").Value = 116.27274
$ws.Range("N97").Value = -2052
$ws.Range("H136").Value = 1994.7858
$ws.Range("I136").Value = 1603.7778
$ws.Range("J136").Value = 2698.6
$ws.Range("K136").Value = 4811.3334
$ws.Range("L136").Value = 8095.799999999999
$ws.Range("M136").Value = -2261.3334
$ws.Range("N136").Value = -13195.8
$ws.Range("H139").Value = 59900
$ws.Range("J139").Value = 59900
$ws.Range("L139").Value = 59900
$ws.Range("N139").Value = -70180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2788.8823
$ws.Range("I105").Value = 2616.5715
$ws.Range("J105").Value = 3593
$ws.Range("K105").Value = 2616.5715
$ws.Range("L105").Value = 3593
$ws.Range("M105").Value = -869.5715
$ws.Range("N105").Value = -7087
$ws.Range("H107").Value = 11859.429
$ws.Range("I107").Value = 11500.5
$ws.Range("K107").Value = 11500.5
$ws.Range("M107").Value = -9580.5
$ws.Range("H134").Value = 3449.6667
$ws.Range("I134").Value = 3306.8096
$ws.Range("K134").Value = 9920.4288
$ws.Range("M134").Value = -7385.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2301.35
$ws.Range("I31").Value = 2277.3333
$ws.Range("J31").Value = 2337.375
$ws.Range("K31").Value = 2277.3333
$ws.Range("L31").Value = 2337.375
$ws.Range("M31").Value = -1982.3333
$ws.Range("N31").Value = -2927.375
$ws.Range("H34").Value = 2301.35
$ws.Range("I34").Value = 2277.3333
$ws.Range("J34").Value = 2337.375
$ws.Range("K34").Value = 2277.3333
$ws.Range("L34").Value = 2337.375
$ws.Range("M34").Value = -2075.3333
$ws.Range("N34").Value = -2741.375
$ws.Range("H105").Value = 23162.125
$ws.Range("J105").Value = 1499.5
$ws.Range("L105").Value = 1499.5
$ws.Range("N105").Value = -4993.5
$ws.Range("H107").Value = 937.8570999999999
$ws.Range("I107").Value = 482.75
$ws.Range("K107").Value = 482.75
$ws.Range("M107").Value = 1437.25
$ws.Range("H132").Value = 2812.9167
$ws.Range("I132").Value = 1630.1666
$ws.Range("J132").Value = 3995.6667
$ws.Range("K132").Value = 4890.4998
$ws.Range("L132").Value = 11987.0001
$ws.Range("M132").Value = -2360.4998
$ws.Range("N132").Value = -17047.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2505.5715
$ws.Range("J5").Value = 3197.8
$ws.Range("L5").Value = 9593.400000000001
$ws.Range("N5").Value = -9817.400000000001
$ws.Range("H7").Value = 669.17645
$ws.Range("I7").Value = 677.1818
$ws.Range("K7").Value = 2031.5454
$ws.Range("M7").Value = -1919.5454
$ws.Range("H35").Value = 900
$ws.Range("J35").Value = 900
$ws.Range("L35").Value = 2700
$ws.Range("N35").Value = -3276
$ws.Range("H36").Value = 216.66667
$ws.Range("I36").Value = 216.66667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 650.00001
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -481.00001
$ws.Range("N36").ClearContents()
$ws.Range("H38").Value = 143.29411
$ws.Range("J38").Value = 170.375
$ws.Range("L38").Value = 511.125
$ws.Range("N38").Value = -1205.125
$ws.Range("H42").Value = 2199
$ws.Range("J42").Value = 2199
$ws.Range("L42").Value = 6597
$ws.Range("N42").Value = -7665
$ws.Range("H69").Value = 7750
$ws.Range("I69").Value = 7750
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 23250
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -22439
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 7750
$ws.Range("I72").Value = 7750
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 69750
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -65694
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 1164
$ws.Range("J132").Value = 1450
$ws.Range("L132").Value = 13050
$ws.Range("N132").Value = -18110
$ws.Range("H135").Value = 2505.5715
$ws.Range("J135").Value = 3197.8
$ws.Range("L135").Value = 28780.2
$ws.Range("N135").Value = -33850.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 53098.9
$ws.Range("I22").Value = 333960
$ws.Range("K22").Value = 333960
$ws.Range("M22").Value = -333665
$ws.Range("H27").Value = 53098.9
$ws.Range("I27").Value = 333960
$ws.Range("K27").Value = 333960
$ws.Range("M27").Value = -333853
$ws.Range("H40").Value = 7157.25
$ws.Range("I40").Value = 6944.4546
$ws.Range("J40").Value = 7417.3335
$ws.Range("K40").Value = 6944.4546
$ws.Range("L40").Value = 7417.3335
$ws.Range("M40").Value = -6808.4546
$ws.Range("N40").Value = -7689.3335
$ws.Range("H43").Value = 19950
$ws.Range("J43").Value = 19950
$ws.Range("L43").Value = 19950
$ws.Range("N43").Value = -20336
$ws.Range("H106").Value = 12705.143
$ws.Range("J106").Value = 12705.143
$ws.Range("L106").Value = 12705.143
$ws.Range("N106").Value = -15229.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2910.6667
$ws.Range("I81").Value = 1934.1818
$ws.Range("J81").Value = 3736.923
$ws.Range("K81").Value = 3868.3636
$ws.Range("L81").Value = 7473.846
$ws.Range("M81").Value = -2807.3636
$ws.Range("N81").Value = -9595.846
$ws.Range("H84").Value = 2910.6667
$ws.Range("I84").Value = 1934.1818
$ws.Range("J84").Value = 3736.923
$ws.Range("K84").Value = 19341.818
$ws.Range("L84").Value = 37369.23
$ws.Range("M84").Value = -14037.818
$ws.Range("N84").Value = -47977.23
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 259.35
$ws.Range("I107").Value = 279.84616
$ws.Range("K107").Value = 839.5384799999999
$ws.Range("M107").Value = 1080.46152

